$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 2.920844475933321
$ws.Range("C6").Value = 3.50332108551396
$ws.Range("D6").Value = 1.643576376937208
$ws.Range("E6").Value = 1.629437769201222
$ws.Range("J6").Value = 2.920844475933321
$ws.Range("K6").Value = 3.50332108551396
$ws.Range("L6").Value = 1.643576376937208
$ws.Range("M6").Value = 1.629437769201222
$ws.Range("B7").Value = 2.62486353443764
$ws.Range("C7").Value = -0.8224280183251671
$ws.Range("D7").Value = 1.589685184540286
$ws.Range("E7").Value = 1.572113170302798
$ws.Range("J7").Value = 2.62486353443764
$ws.Range("K7").Value = -0.8224280183251671
$ws.Range("L7").Value = 1.589685184540286
$ws.Range("M7").Value = 1.572113170302798
$ws.Range("B8").Value = 2.478507513310722
$ws.Range("C8").Value = 0.6973179972408781
$ws.Range("D8").Value = 1.182006459958566
$ws.Range("E8").Value = 1.491313139750291
$ws.Range("J8").Value = 2.478507513310722
$ws.Range("K8").Value = 0.6973179972408781
$ws.Range("L8").Value = 1.182006459958566
$ws.Range("M8").Value = 1.491313139750291
$ws.Range("B9").Value = -0.5851866069045748
$ws.Range("C9").Value = 3.859827423951487
$ws.Range("D9").Value = 1.549208067845804
$ws.Range("E9").Value = 1.601297700235178
$ws.Range("J9").Value = -0.5851866069045748
$ws.Range("K9").Value = 3.859827423951487
$ws.Range("L9").Value = 1.549208067845804
$ws.Range("M9").Value = 1.601297700235178
$ws.Range("B10").Value = 1.532569976950621
$ws.Range("C10").Value = 2.359928885994272
$ws.Range("D10").Value = 1.930538557116238
$ws.Range("E10").Value = 1.654336769972793
$ws.Range("J10").Value = 1.532569976950621
$ws.Range("K10").Value = 2.359928885994272
$ws.Range("L10").Value = 1.930538557116238
$ws.Range("M10").Value = 1.654336769972793
$ws.Range("B11").Value = 4.49236174839908
$ws.Range("C11").Value = 5.619242251286392
$ws.Range("D11").Value = 2.045174124597719
$ws.Range("E11").Value = 1.726261054191536
$ws.Range("J11").Value = 4.49236174839908
$ws.Range("K11").Value = 5.619242251286392
$ws.Range("L11").Value = 2.045174124597719
$ws.Range("M11").Value = 1.726261054191536
$ws.Range("B12").Value = 2.426354921777689
$ws.Range("C12").Value = 6.511403028606294
$ws.Range("D12").Value = 2.565404917870143
$ws.Range("E12").Value = 1.883318937686796
$ws.Range("J12").Value = 2.426354921777689
$ws.Range("K12").Value = 6.511403028606294
$ws.Range("L12").Value = 2.565404917870143
$ws.Range("M12").Value = 1.883318937686796
$ws.Range("B13").Value = 3.453365118268367
$ws.Range("C13").Value = 3.574430157754606
$ws.Range("D13").Value = 2.617699030388895
$ws.Range("E13").Value = 1.899968678029008
$ws.Range("J13").Value = 3.453365118268367
$ws.Range("K13").Value = 3.574430157754606
$ws.Range("L13").Value = 2.617699030388895
$ws.Range("M13").Value = 1.899968678029008
$ws.Range("B14").Value = 4.457383385483923
$ws.Range("C14").Value = 7.812108007885445
$ws.Range("D14").Value = 2.603983628082156
$ws.Range("E14").Value = 2.003499679048085
$ws.Range("J14").Value = 4.457383385483923
$ws.Range("K14").Value = 7.812108007885445
$ws.Range("L14").Value = 2.603983628082156
$ws.Range("M14").Value = 2.003499679048085
$ws.Range("B15").Value = 4.741527240123398
$ws.Range("C15").Value = 8.035795149322816
$ws.Range("D15").Value = 3.246591233005434
$ws.Range("E15").Value = 2.200677940063331
$ws.Range("J15").Value = 4.741527240123398
$ws.Range("K15").Value = 8.035795149322816
$ws.Range("L15").Value = 3.246591233005434
$ws.Range("M15").Value = 2.200677940063331
$ws.Range("B16").Value = 4.914510101920945
$ws.Range("C16").Value = 8.462360184962892
$ws.Range("D16").Value = 3.518084194922744
$ws.Range("E16").Value = 2.324739161053456
$ws.Range("J16").Value = 4.914510101920945
$ws.Range("K16").Value = 8.462360184962892
$ws.Range("L16").Value = 3.518084194922744
$ws.Range("M16").Value = 2.324739161053456
$ws.Range("B17").Value = 4.530987644907642
$ws.Range("C17").Value = 1.086403302456671
$ws.Range("D17").Value = 3.159679259395985
$ws.Range("E17").Value = 2.228979836847881
$ws.Range("J17").Value = 4.530987644907642
$ws.Range("K17").Value = 1.086403302456671
$ws.Range("L17").Value = 3.159679259395985
$ws.Range("M17").Value = 2.228979836847881
$ws.Range("B18").Value = 4.152352593920312
$ws.Range("C18").Value = 3.829086583649435
$ws.Range("D18").Value = 2.308001588488731
$ws.Range("E18").Value = 2.173252283259402
$ws.Range("J18").Value = 4.152352593920312
$ws.Range("K18").Value = 3.829086583649435
$ws.Range("L18").Value = 2.308001588488731
$ws.Range("M18").Value = 2.173252283259402
$ws.Range("B19").Value = 4.17161561064515
$ws.Range("C19").Value = 3.426370195286472
$ws.Range("D19").Value = 2.578060717518183
$ws.Range("E19").Value = 2.255472600394167
$ws.Range("J19").Value = 4.17161561064515
$ws.Range("K19").Value = 3.426370195286472
$ws.Range("L19").Value = 2.578060717518183
$ws.Range("M19").Value = 2.255472600394167
